# Update Naming Plots Regression
# Clear the leftover regression-parameter columns (AL, AR:AZ, BB, BC, BE, BF, BH)
# for rows 18 through 32 on the "PM3_Pu_GHS" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PM3_Pu_GHS")

for ($r = 18; $r -le 32; $r++) {
    $ws.Range("AL${r}").ClearContents()
    $ws.Range("AR${r}:AZ${r}").ClearContents()
    $ws.Range("BB${r}:BC${r}").ClearContents()
    $ws.Range("BE${r}:BF${r}").ClearContents()
    $ws.Range("BH${r}").ClearContents()
}
